# Weekly price-sheet update: a new week's record for "Berenjena" at
# "Terminal Hortofrutícola Agro Chillán" is inserted at the top of the
# data block (row 95), pushing every existing record down by one row
# (old row 125 becomes row 126). This mirrors how the source feed
# prepends the latest week and keeps the historical rows intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 95; rows 95:125 shift down to 96:126.
$ws.Rows("95:95").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A95").Value = 7
$ws.Range("B95").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C95").Value = "Ñuble"
$ws.Range("D95").Value = 45211
$ws.Range("E95").Value = 16
$ws.Range("F95").Value = 100112001
$ws.Range("G95").Value = "Berenjena"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 50
$ws.Range("K95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("M95").Value = 10000
$ws.Range("N95").Value = "$/caja 60 unidades"
$ws.Range("O95").Value = "Región de Arica y Parinacota"
$ws.Range("P95").Value = 167
$ws.Range("Q95").Value = 60
$ws.Range("R95").Value = "Hortaliza"
